# "revisi menrevisi kombinasi loc"
# The LOC-combination results table (A:G, header-less) is recomputed:
#  - a new row is inserted for FY_4.png (old row 36 "TO_1.png" ... shift down to row 37),
#    pushing the sheet dimension from A1:G44 to A1:G45
#  - every row's distance/probability metrics (columns B:E) are updated
#  - several rows' recognised-name (F) / verdict (G) cells flip as a result
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "AK_1.png"
$ws.Cells.Item(1, 2).Value = 0.472
$ws.Cells.Item(1, 3).Value = 0.002
$ws.Cells.Item(1, 4).Value = 0.605
$ws.Cells.Item(1, 5).Value = 0.3
$ws.Cells.Item(1, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(1, 7).Value = "Salah"

$ws.Cells.Item(2, 1).Value = "AK_2.png"
$ws.Cells.Item(2, 2).Value = 0.861
$ws.Cells.Item(2, 3).Value = 0.003
$ws.Cells.Item(2, 4).Value = 0.823
$ws.Cells.Item(2, 5).Value = 0.533
$ws.Cells.Item(2, 6).Value = "Akhlak Kamiswara"
$ws.Cells.Item(2, 7).Value = "Benar"

$ws.Cells.Item(3, 1).Value = "AK_3.png"
$ws.Cells.Item(3, 2).Value = 0.842
$ws.Cells.Item(3, 3).Value = 0.003
$ws.Cells.Item(3, 4).Value = 0.825
$ws.Cells.Item(3, 5).Value = 0.5669999999999999
$ws.Cells.Item(3, 6).Value = "Akhlak Kamiswara"
$ws.Cells.Item(3, 7).Value = "Benar"

$ws.Cells.Item(4, 1).Value = "AK_4.png"
$ws.Cells.Item(4, 2).Value = 0.435
$ws.Cells.Item(4, 3).Value = 0.001
$ws.Cells.Item(4, 4).Value = 0.593
$ws.Cells.Item(4, 5).Value = 0.3
$ws.Cells.Item(4, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(4, 7).Value = "Salah"

$ws.Cells.Item(5, 1).Value = "AK_5.png"
$ws.Cells.Item(5, 2).Value = 0.442
$ws.Cells.Item(5, 3).Value = 0.001
$ws.Cells.Item(5, 4).Value = 0.83
$ws.Cells.Item(5, 5).Value = 0.333
$ws.Cells.Item(5, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(5, 7).Value = "Salah"

$ws.Cells.Item(6, 1).Value = "MIB_1.png"
$ws.Cells.Item(6, 2).Value = 1.305
$ws.Cells.Item(6, 3).Value = 0.004
$ws.Cells.Item(6, 4).Value = 0.575
$ws.Cells.Item(6, 5).Value = 0.667
$ws.Cells.Item(6, 6).Value = "Muhammad Iqbal Baqi"
$ws.Cells.Item(6, 7).Value = "Benar"

$ws.Cells.Item(7, 1).Value = "MIB_2.png"
$ws.Cells.Item(7, 2).Value = 1.098
$ws.Cells.Item(7, 3).Value = 0.004
$ws.Cells.Item(7, 4).Value = 0.732
$ws.Cells.Item(7, 5).Value = 0.7
$ws.Cells.Item(7, 6).Value = "Muhammad Iqbal Baqi"
$ws.Cells.Item(7, 7).Value = "Benar"

$ws.Cells.Item(8, 1).Value = "MIB_3.png"
$ws.Cells.Item(8, 2).Value = 1.315
$ws.Cells.Item(8, 3).Value = 0.004
$ws.Cells.Item(8, 4).Value = 0.909
$ws.Cells.Item(8, 5).Value = 0.967
$ws.Cells.Item(8, 6).Value = "Muhammad Iqbal Baqi"
$ws.Cells.Item(8, 7).Value = "Benar"

$ws.Cells.Item(9, 1).Value = "MIB_4.png"
$ws.Cells.Item(9, 2).Value = 1.1
$ws.Cells.Item(9, 3).Value = 0.004
$ws.Cells.Item(9, 4).Value = 0.705
$ws.Cells.Item(9, 5).Value = 0.7
$ws.Cells.Item(9, 6).Value = "Muhammad Iqbal Baqi"
$ws.Cells.Item(9, 7).Value = "Benar"

$ws.Cells.Item(10, 1).Value = "MIB_5.png"
$ws.Cells.Item(10, 2).Value = 1.271
$ws.Cells.Item(10, 3).Value = 0.004
$ws.Cells.Item(10, 4).Value = 0.718
$ws.Cells.Item(10, 5).Value = 0.733
$ws.Cells.Item(10, 6).Value = "Muhammad Iqbal Baqi"
$ws.Cells.Item(10, 7).Value = "Benar"

$ws.Cells.Item(11, 1).Value = "AAH_1.png"
$ws.Cells.Item(11, 2).Value = 0.824
$ws.Cells.Item(11, 3).Value = 0.003
$ws.Cells.Item(11, 4).Value = 0.75
$ws.Cells.Item(11, 5).Value = 0.9
$ws.Cells.Item(11, 6).Value = "Andrea Ayunove Hutami"
$ws.Cells.Item(11, 7).Value = "Benar"

$ws.Cells.Item(12, 1).Value = "AAH_2.png"
$ws.Cells.Item(12, 2).Value = 1.124
$ws.Cells.Item(12, 3).Value = 0.004
$ws.Cells.Item(12, 4).Value = 0.91
$ws.Cells.Item(12, 5).Value = 0.967
$ws.Cells.Item(12, 6).Value = "Andrea Ayunove Hutami"
$ws.Cells.Item(12, 7).Value = "Benar"

$ws.Cells.Item(13, 1).Value = "AAH_3.png"
$ws.Cells.Item(13, 2).Value = 0.84
$ws.Cells.Item(13, 3).Value = 0.003
$ws.Cells.Item(13, 4).Value = 0.759
$ws.Cells.Item(13, 5).Value = 0.967
$ws.Cells.Item(13, 6).Value = "Andrea Ayunove Hutami"
$ws.Cells.Item(13, 7).Value = "Benar"

$ws.Cells.Item(14, 1).Value = "TI_1.png"
$ws.Cells.Item(14, 2).Value = 0.975
$ws.Cells.Item(14, 3).Value = 0.003
$ws.Cells.Item(14, 4).Value = 0.792
$ws.Cells.Item(14, 5).Value = 0.3
$ws.Cells.Item(14, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(14, 7).Value = "Salah"

$ws.Cells.Item(15, 1).Value = "TI_2.png"
$ws.Cells.Item(15, 2).Value = 0.9379999999999999
$ws.Cells.Item(15, 3).Value = 0.003
$ws.Cells.Item(15, 4).Value = 0.805
$ws.Cells.Item(15, 5).Value = 0.5
$ws.Cells.Item(15, 6).Value = "Toni Ismail"
$ws.Cells.Item(15, 7).Value = "Benar"

$ws.Cells.Item(16, 1).Value = "TI_3.png"
$ws.Cells.Item(16, 2).Value = 0.717
$ws.Cells.Item(16, 3).Value = 0.002
$ws.Cells.Item(16, 4).Value = 0.901
$ws.Cells.Item(16, 5).Value = 0.333
$ws.Cells.Item(16, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(16, 7).Value = "Salah"

$ws.Cells.Item(17, 1).Value = "TI_4.png"
$ws.Cells.Item(17, 2).Value = 0.6870000000000001
$ws.Cells.Item(17, 3).Value = 0.002
$ws.Cells.Item(17, 4).Value = 0.764
$ws.Cells.Item(17, 5).Value = 0.533
$ws.Cells.Item(17, 6).Value = "Toni Ismail"
$ws.Cells.Item(17, 7).Value = "Benar"

$ws.Cells.Item(18, 1).Value = "TI_5.png"
$ws.Cells.Item(18, 2).Value = 1.011
$ws.Cells.Item(18, 3).Value = 0.003
$ws.Cells.Item(18, 4).Value = 0.854
$ws.Cells.Item(18, 5).Value = 0.633
$ws.Cells.Item(18, 6).Value = "Toni Ismail"
$ws.Cells.Item(18, 7).Value = "Benar"

$ws.Cells.Item(19, 1).Value = "RAS_1.png"
$ws.Cells.Item(19, 2).Value = 0.646
$ws.Cells.Item(19, 3).Value = 0.002
$ws.Cells.Item(19, 4).Value = 0.794
$ws.Cells.Item(19, 5).Value = 0.267
$ws.Cells.Item(19, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(19, 7).Value = "Salah"

$ws.Cells.Item(20, 1).Value = "RAS_2.png"
$ws.Cells.Item(20, 2).Value = 1.106
$ws.Cells.Item(20, 3).Value = 0.004
$ws.Cells.Item(20, 4).Value = 0.75
$ws.Cells.Item(20, 5).Value = 0.467
$ws.Cells.Item(20, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(20, 7).Value = "Salah"

$ws.Cells.Item(21, 1).Value = "RAS_3.png"
$ws.Cells.Item(21, 2).Value = 0.513
$ws.Cells.Item(21, 3).Value = 0.002
$ws.Cells.Item(21, 4).Value = 0.665
$ws.Cells.Item(21, 5).Value = 0.367
$ws.Cells.Item(21, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(21, 7).Value = "Salah"

$ws.Cells.Item(22, 1).Value = "RAS_4.png"
$ws.Cells.Item(22, 2).Value = 1.16
$ws.Cells.Item(22, 3).Value = 0.004
$ws.Cells.Item(22, 4).Value = 0.607
$ws.Cells.Item(22, 5).Value = 0.4
$ws.Cells.Item(22, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(22, 7).Value = "Salah"

$ws.Cells.Item(23, 1).Value = "RAS_5.png"
$ws.Cells.Item(23, 2).Value = 1.078
$ws.Cells.Item(23, 3).Value = 0.004
$ws.Cells.Item(23, 4).Value = 0.784
$ws.Cells.Item(23, 5).Value = 0.533
$ws.Cells.Item(23, 6).Value = "Ridha Ayu Salsabila"
$ws.Cells.Item(23, 7).Value = "Benar"

$ws.Cells.Item(24, 1).Value = "RR_1.png"
$ws.Cells.Item(24, 2).Value = 1.148
$ws.Cells.Item(24, 3).Value = 0.004
$ws.Cells.Item(24, 4).Value = 0.875
$ws.Cells.Item(24, 5).Value = 0.5
$ws.Cells.Item(24, 6).Value = "Rafiqo Rapitasari"
$ws.Cells.Item(24, 7).Value = "Benar"

$ws.Cells.Item(25, 1).Value = "RR_2.png"
$ws.Cells.Item(25, 2).Value = 1.225
$ws.Cells.Item(25, 3).Value = 0.004
$ws.Cells.Item(25, 4).Value = 0.876
$ws.Cells.Item(25, 5).Value = 0.5669999999999999
$ws.Cells.Item(25, 6).Value = "Rafiqo Rapitasari"
$ws.Cells.Item(25, 7).Value = "Benar"

$ws.Cells.Item(26, 1).Value = "RR_3.png"
$ws.Cells.Item(26, 2).Value = 0.9419999999999999
$ws.Cells.Item(26, 3).Value = 0.003
$ws.Cells.Item(26, 4).Value = 0.48
$ws.Cells.Item(26, 5).Value = 0.5
$ws.Cells.Item(26, 6).Value = "Rafiqo Rapitasari"
$ws.Cells.Item(26, 7).Value = "Benar"

$ws.Cells.Item(27, 1).Value = "RR_4.png"
$ws.Cells.Item(27, 2).Value = 1.23
$ws.Cells.Item(27, 3).Value = 0.004
$ws.Cells.Item(27, 4).Value = 0.886
$ws.Cells.Item(27, 5).Value = 0.5
$ws.Cells.Item(27, 6).Value = "Rafiqo Rapitasari"
$ws.Cells.Item(27, 7).Value = "Benar"

$ws.Cells.Item(28, 1).Value = "RR_5.png"
$ws.Cells.Item(28, 2).Value = 1.185
$ws.Cells.Item(28, 3).Value = 0.004
$ws.Cells.Item(28, 4).Value = 0.882
$ws.Cells.Item(28, 5).Value = 0.533
$ws.Cells.Item(28, 6).Value = "Rafiqo Rapitasari"
$ws.Cells.Item(28, 7).Value = "Benar"

$ws.Cells.Item(29, 1).Value = "AR_1.png"
$ws.Cells.Item(29, 2).Value = 0.854
$ws.Cells.Item(29, 3).Value = 0.003
$ws.Cells.Item(29, 4).Value = 0.746
$ws.Cells.Item(29, 5).Value = 0.667
$ws.Cells.Item(29, 6).Value = "Arizli Romadhon"
$ws.Cells.Item(29, 7).Value = "Benar"

$ws.Cells.Item(30, 1).Value = "GA_1.png"
$ws.Cells.Item(30, 2).Value = 1.189
$ws.Cells.Item(30, 3).Value = 0.004
$ws.Cells.Item(30, 4).Value = 0.804
$ws.Cells.Item(30, 5).Value = 0.9330000000000001
$ws.Cells.Item(30, 6).Value = "Gege Ardiyansyah"
$ws.Cells.Item(30, 7).Value = "Benar"

$ws.Cells.Item(31, 1).Value = "GA_2.png"
$ws.Cells.Item(31, 2).Value = 0.673
$ws.Cells.Item(31, 3).Value = 0.002
$ws.Cells.Item(31, 4).Value = 0.773
$ws.Cells.Item(31, 5).Value = 0.233
$ws.Cells.Item(31, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(31, 7).Value = "Salah"

$ws.Cells.Item(32, 1).Value = "GA_3.png"
$ws.Cells.Item(32, 2).Value = 0.733
$ws.Cells.Item(32, 3).Value = 0.002
$ws.Cells.Item(32, 4).Value = 0.728
$ws.Cells.Item(32, 5).Value = 0.3
$ws.Cells.Item(32, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(32, 7).Value = "Salah"

$ws.Cells.Item(33, 1).Value = "FY_1.png"
$ws.Cells.Item(33, 2).Value = 1.079
$ws.Cells.Item(33, 3).Value = 0.004
$ws.Cells.Item(33, 4).Value = 0.766
$ws.Cells.Item(33, 5).Value = 0.233
$ws.Cells.Item(33, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(33, 7).Value = "Salah"

$ws.Cells.Item(34, 1).Value = "FY_2.png"
$ws.Cells.Item(34, 2).Value = 1.283
$ws.Cells.Item(34, 3).Value = 0.004
$ws.Cells.Item(34, 4).Value = 0.724
$ws.Cells.Item(34, 5).Value = 0.5
$ws.Cells.Item(34, 6).Value = "Andrea Ayunove Hutami"
$ws.Cells.Item(34, 7).Value = "Salah"

$ws.Cells.Item(35, 1).Value = "FY_3.png"
$ws.Cells.Item(35, 2).Value = 1.239
$ws.Cells.Item(35, 3).Value = 0.004
$ws.Cells.Item(35, 4).Value = 0.781
$ws.Cells.Item(35, 5).Value = 0.367
$ws.Cells.Item(35, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(35, 7).Value = "Salah"

$ws.Cells.Item(36, 1).Value = "FY_4.png"
$ws.Cells.Item(36, 2).Value = 1.252
$ws.Cells.Item(36, 3).Value = 0.004
$ws.Cells.Item(36, 4).Value = 0.8080000000000001
$ws.Cells.Item(36, 5).Value = 0.233
$ws.Cells.Item(36, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(36, 7).Value = "Salah"

$ws.Cells.Item(37, 1).Value = "TO_1.png"
$ws.Cells.Item(37, 2).Value = 0.802
$ws.Cells.Item(37, 3).Value = 0.003
$ws.Cells.Item(37, 4).Value = 0.772
$ws.Cells.Item(37, 5).Value = 0.533
$ws.Cells.Item(37, 6).Value = "Tiara Oktavian"
$ws.Cells.Item(37, 7).Value = "Benar"

$ws.Cells.Item(38, 1).Value = "TO_2.png"
$ws.Cells.Item(38, 2).Value = 0.92
$ws.Cells.Item(38, 3).Value = 0.003
$ws.Cells.Item(38, 4).Value = 0.824
$ws.Cells.Item(38, 5).Value = 0.533
$ws.Cells.Item(38, 6).Value = "Tiara Oktavian"
$ws.Cells.Item(38, 7).Value = "Benar"

$ws.Cells.Item(39, 1).Value = "TO_3.png"
$ws.Cells.Item(39, 2).Value = 0.825
$ws.Cells.Item(39, 3).Value = 0.003
$ws.Cells.Item(39, 4).Value = 0.8100000000000001
$ws.Cells.Item(39, 5).Value = 0.533
$ws.Cells.Item(39, 6).Value = "Tiara Oktavian"
$ws.Cells.Item(39, 7).Value = "Benar"

$ws.Cells.Item(40, 1).Value = "TO_4.png"
$ws.Cells.Item(40, 2).Value = 2.203
$ws.Cells.Item(40, 3).Value = 0.007
$ws.Cells.Item(40, 4).Value = 0.384
$ws.Cells.Item(40, 5).Value = 0.4
$ws.Cells.Item(40, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(40, 7).Value = "Salah"

$ws.Cells.Item(41, 1).Value = "TO_5.png"
$ws.Cells.Item(41, 2).Value = 1.953
$ws.Cells.Item(41, 3).Value = 0.006
$ws.Cells.Item(41, 4).Value = 0.251
$ws.Cells.Item(41, 5).Value = 0.333
$ws.Cells.Item(41, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(41, 7).Value = "Salah"

$ws.Cells.Item(42, 1).Value = "TD_1.png"
$ws.Cells.Item(42, 2).Value = 1.481
$ws.Cells.Item(42, 3).Value = 0.005
$ws.Cells.Item(42, 4).Value = 0.374
$ws.Cells.Item(42, 5).Value = 0.367
$ws.Cells.Item(42, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(42, 7).Value = "Benar"

$ws.Cells.Item(43, 1).Value = "TD_2.png"
$ws.Cells.Item(43, 2).Value = 1.479
$ws.Cells.Item(43, 3).Value = 0.005
$ws.Cells.Item(43, 4).Value = 0.345
$ws.Cells.Item(43, 5).Value = 0.3
$ws.Cells.Item(43, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(43, 7).Value = "Benar"

$ws.Cells.Item(44, 1).Value = "TD_3.png"
$ws.Cells.Item(44, 2).Value = 1.018
$ws.Cells.Item(44, 3).Value = 0.003
$ws.Cells.Item(44, 4).Value = 0.703
$ws.Cells.Item(44, 5).Value = 0.233
$ws.Cells.Item(44, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(44, 7).Value = "Benar"

$ws.Cells.Item(45, 1).Value = "TD_4.png"
$ws.Cells.Item(45, 2).Value = 0.97
$ws.Cells.Item(45, 3).Value = 0.003
$ws.Cells.Item(45, 4).Value = 0.665
$ws.Cells.Item(45, 5).Value = 0.267
$ws.Cells.Item(45, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(45, 7).Value = "Benar"
